$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Suppress print statements?" bullet paragraph entirely
#    (it sat between "Feed win or loss back to AI for learning" and
#    "Make code cleaner").
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Suppress print statements?")
if ($found) {
    $paraRange = $rng.Paragraphs(1).Range
    $paraRange.Delete()
}

# ---------------------------------------------------------------------------
# 2. Mark the built-in "Default Paragraph Font" style as semi-hidden, i.e.
#    not offered/recommended in the UI (w:semiHidden), matching the
#    Visibility = False state Word uses for that flag.
# ---------------------------------------------------------------------------
$style = $d.Styles("Default Paragraph Font")
$style.Visibility = $false
